# Auto: Weekly update of data
# Appends 8 new match rows (107-114) to the "Main" sheet of the CS2 pre-match
# statistics workbook, then restores the frozen-header view / selection state
# so the sheet looks the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$colIndex = @{ 'A'=1; 'B'=2; 'C'=3; 'D'=4; 'E'=5; 'F'=6; 'G'=7; 'H'=8; 'I'=9; 'J'=10; 'K'=11; 'L'=12; 'M'=13; 'N'=14; 'O'=15; 'P'=16; 'Q'=17; 'R'=18; 'S'=19; 'T'=20 }

$newRows = @{
    107 = @{ 'A'=866; 'B'=2111; 'C'=221; 'D'=18164; 'E'=18082; 'F'=19653; 'H'=17735; 'I'=18315; 'J'=21268; 'K'=13810; 'L'=14328; 'M'=19893; 'N'=5; 'O'=5; 'P'=0; 'Q'=-3; 'R'=-3; 'S'=-3; 'T'=-1 }
    108 = @{ 'A'=867; 'B'=2111; 'C'=5; 'D'=16887; 'E'=20188; 'F'=18478; 'G'=18316; 'H'=17717; 'I'=12205; 'J'=11274; 'L'=15942; 'M'=22509; 'N'=1; 'O'=6; 'P'=4; 'Q'=4; 'R'=-5; 'S'=-12; 'T'=2 }
    109 = @{ 'A'=868; 'B'=5; 'C'=2111; 'D'=25254; 'E'=24694; 'F'=20798; 'G'=17819; 'H'=23525; 'I'=22813; 'K'=22349; 'L'=20112; 'M'=23701; 'N'=-3; 'O'=1; 'P'=0; 'Q'=8; 'R'=1; 'S'=3; 'T'=-3 }
    110 = @{ 'A'=869; 'B'=5; 'C'=41; 'D'=25595; 'E'=24999; 'F'=23838; 'G'=18126; 'H'=21151; 'I'=25102; 'J'=24999; 'K'=20875; 'L'=22325; 'M'=22424; 'N'=-11; 'O'=11; 'P'=0; 'Q'=-3; 'R'=6; 'S'=6; 'T'=-9 }
    111 = @{ 'A'=870; 'B'=311; 'C'=221; 'D'=18146; 'E'=19999; 'F'=18683; 'G'=14107; 'H'=17837; 'I'=18175; 'J'=17368; 'K'=18116; 'L'=17940; 'M'=17551; 'N'=1; 'O'=-2; 'P'=4; 'Q'=-2; 'R'=3; 'S'=5; 'T'=-9 }
    112 = @{ 'A'=871; 'B'=41; 'C'=41; 'D'=21649; 'E'=21017; 'F'=23606; 'G'=18196; 'H'=21558; 'I'=21628; 'J'=19618; 'K'=23036; 'M'=20001; 'N'=1; 'O'=-2; 'P'=4; 'Q'=-1; 'R'=2; 'S'=0; 'T'=-4 }
    113 = @{ 'A'=872; 'B'=41; 'C'=221; 'D'=22779; 'E'=21998; 'F'=21340; 'G'=24100; 'H'=18542; 'I'=23289; 'J'=24777; 'K'=25842; 'L'=20000; 'M'=22300; 'N'=-6; 'O'=4; 'P'=2; 'Q'=-11; 'R'=1; 'S'=4; 'T'=6 }
    114 = @{ 'A'=873; 'B'=5; 'C'=2111; 'D'=21696; 'E'=18893; 'F'=22372; 'G'=17122; 'H'=24628; 'I'=22122; 'J'=23365; 'L'=21733; 'M'=21457; 'N'=2; 'O'=7; 'P'=4; 'Q'=-6; 'R'=-2; 'S'=5; 'T'=-3 }
}

foreach ($r in ($newRows.Keys | Sort-Object)) {
    $rowData = $newRows[$r]
    foreach ($colLetter in $rowData.Keys) {
        $colNum = $colIndex[$colLetter]
        $ws.Cells.Item([int]$r, $colNum).Value2 = $rowData[$colLetter]
    }
}

# Restore the view: header row frozen, scrolled so row 98 is the first
# visible row under the freeze, with the last-entered cell selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 98
$ws.Range("U114").Select()

